# Update "想去人数" (F column) figures across the three affected sheets.
# Values were refreshed from the live data source (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 782
$ws1.Range("F3").Value = 59
$ws1.Range("F4").Value = 412
$ws1.Range("F5").Value = 140
$ws1.Range("F7").Value = 157
$ws1.Range("F8").Value = 343
$ws1.Range("F9").Value = 453
$ws1.Range("F10").Value = 511
$ws1.Range("F12").Value = 11755
$ws1.Range("F13").Value = 5412

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 106
$ws2.Range("F3").Value = 3

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 782
$ws4.Range("F3").Value = 59
$ws4.Range("F4").Value = 106
$ws4.Range("F5").Value = 3
$ws4.Range("F6").Value = 412
$ws4.Range("F7").Value = 140
$ws4.Range("F9").Value = 157
$ws4.Range("F10").Value = 343
$ws4.Range("F11").Value = 453
$ws4.Range("F12").Value = 511
$ws4.Range("F14").Value = 11755
$ws4.Range("F16").Value = 5412
